# Apply stimulus-set rename (face -> book) and correct_ans code expansion
# (b/r/y -> center/right/left) across the main_conditions_23 sequence sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$ansMap = @{ "b" = "center"; "r" = "right"; "y" = "left" }

for ($r = 2; $r -le $lastRow; $r++) {

    # Columns A-D hold the "<category>//<category>_NN.jpg" filenames used
    # for this trial's prompt/correct/distractor images. Any reference to
    # the old "face" stimulus set now points at the new "book" set.
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.StartsWith("face//face")) {
            $cell.Value = $val.Replace("face//face", "book//book")
        }
    }

    # Column L ("correct_ans") used cryptic single-letter codes; expand them
    # to the full, human-readable position names.
    $lCell = $ws.Cells.Item($r, 12)
    $lVal = $lCell.Value2
    if ($ansMap.ContainsKey($lVal)) {
        $lCell.Value = $ansMap[$lVal]
    }
}
